$d = $word.ActiveDocument

# --- Edit 1: expand the intro paragraph describing individual contributions to the
#             tidy data sets, and mention the one group member's template work. ---
$oldIntro = "For the final project, our group was tasked with the input and output of the data via the Shiny app. This involved appropriately `"tidying`" the different data tables and outputting the other groups' graphs and plots in the Shiny app. Since there were four different data tables, each group member was challenged to tidy one of the data sets. The data sets were labeled: efficacy, plasma, tissue laser, and tissue standard pk. Each member of our group describes the work they did individually. This is followed by challenges faced and lesson learned by the group as a whole."
$newIntro = "For the final project, our group was tasked with the input and output of the data via the Shiny app. This involved appropriately `"tidying`" the different data tables and outputting the other groups' graphs and plots in the Shiny app. Since there were four different data tables, each group member was challenged to tidy one of the data sets. The data sets were labeled: efficacy, plasma, tissue laser, and tissue standard pk. Each member of our group describes the work they did individually to those data sets. One group also describes the work she did to an original template of the data we received. These descriptions are followed by challenges faced and lesson learned by the group as a whole."

$introPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq $oldIntro) {
        $introPara = $d.Paragraphs($i)
        break
    }
}
if ($introPara -eq $null) {
    throw "Could not locate the intro paragraph to edit"
}
$introPara.Range.Text = $newIntro

# --- Edit 2: append a new closing "Overall" paragraph (Body Text style) after the
#             paragraph about individual mice / measurements, just before the end
#             of the document. ---
$mousePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("It was interesting to think about the challenges surrounding identifying individual mice")) {
        $mousePara = $d.Paragraphs($i)
    }
}
if ($mousePara -eq $null) {
    throw "Could not locate the closing mouse-data paragraph"
}

$insertionPoint = $d.Range($mousePara.Range.End, $mousePara.Range.End)
$insertionPoint.Text = "`rOverall, we learned a lot about the importance of `"tidy`" data using real data. It was an important lesson to learn that not all data we may receive will be clean. However, it was rewarding to take the research group's data collected from a lot of time and hard work and be able to give them a Shiny app that allows them to visually and analytically explore their data all in one go."

$newPara = $d.Paragraphs.Last
$newPara.Style = "BodyText"

Write-Output "edits applied"
